$d = $word.ActiveDocument

# 1) Job Designation: Graphics Designer -> Software Developer
$d.Content.Find.Execute("Graphics Designer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Software Developer", 2) | Out-Null

# 2) Place of Posting (both occurrences): Hyderabad -> Delhi
$d.Content.Find.Execute("Hyderabad", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Delhi", 2) | Out-Null

# 3) Admin Information contact table (Table 4): fill in the previously-blank
#    "Second Contact Person" cells.
$contactTable = $d.Tables.Item(4)
$contactTable.Cell(2, 3).Range.Text = "Rahul Gupta"     # Name row
$contactTable.Cell(4, 3).Range.Text = "1234567890"      # Email Address row
$contactTable.Cell(5, 3).Range.Text = "1234567890"      # Mobile Number row

# 4) Placement calendar / mode-of-hiring table (Table 6)
$stagesTable = $d.Tables.Item(6)
$stagesTable.Cell(5, 3).Range.Text = "Virtual"         # Group Discussion
$stagesTable.Cell(7, 3).Range.Text = "Not Applicable"  # Any other rounds
